$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E22").Value = "accomodation"
$ws.Range("F22").Value = "accommodation"
$ws.Range("E23").Value = "enourmous"
$ws.Range("F23").Value = "enormous"
$ws.Range("E24").Value = "unevitably"
$ws.Range("F24").Value = "inevitably"

$excel.ActiveWindow.ScrollRow = 3
$ws.Range("I21").Select()
